# Fix a couple bugs: correct the values in column D (third numeric column)
# on Sheet1 which had been mistakenly duplicated from column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value  = 1901
$ws.Range("D5").Value  = 749
$ws.Range("D6").Value  = 159
$ws.Range("D7").Value  = 85
$ws.Range("D8").Value  = 162
$ws.Range("D10").Value = 49
$ws.Range("D13").Value = 1797
$ws.Range("D14").Value = 406
$ws.Range("D17").Value = 2.098
$ws.Range("D20").Value = 0
$ws.Range("D23").Value = 6176
$ws.Range("D26").Value = 296
$ws.Range("D28").Value = 437
$ws.Range("D29").Value = 574
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 312
$ws.Range("D32").Value = 3965
$ws.Range("D33").Value = 4.277
$ws.Range("D35").Value = 249
$ws.Range("D37").Value = 35
